$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the alias-name row (row 1, the table header / display name) with the
# field-name row (row 3) so the table header now shows the Chinese display
# names and row 3 holds the original English field names.
$row1 = @("序列", "比赛id", "偏移", "开始日期", "左玩家类型", "左玩家id", "右玩家类型", "右玩家id")
$row3 = @("Id", "Tid", "Offset", "Date", "LeftType", "LeftValue", "RightType", "RightValue")

for ($i = 0; $i -lt 8; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $row1[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

# Move the active selection to D12, matching the saved view state.
$ws.Range("D12").Select()
